$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header: add "time_taken" column in F1, copying header style from E1
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F1").Value = "time_taken"

# Timestamps for rows 2-299 (time_taken metadata column)
$timeTaken = @(
    "2021-10-05 13:40:13.764511",
    "2021-10-05 13:40:13.764522",
    "2021-10-05 13:40:13.764526",
    "2021-10-05 13:40:13.764529",
    "2021-10-05 13:40:13.764531",
    "2021-10-05 13:40:13.764534",
    "2021-10-05 13:40:13.764537",
    "2021-10-05 13:40:13.764539",
    "2021-10-05 13:40:13.764542",
    "2021-10-05 13:40:13.764545",
    "2021-10-05 13:40:13.764547",
    "2021-10-05 13:40:13.764550",
    "2021-10-05 13:40:13.764553",
    "2021-10-05 13:40:13.764555",
    "2021-10-05 13:40:13.764558",
    "2021-10-05 13:40:13.764560",
    "2021-10-05 13:40:13.764563",
    "2021-10-05 13:40:13.764566",
    "2021-10-05 13:40:13.764569",
    "2021-10-05 13:40:13.764571",
    "2021-10-05 13:40:13.764574",
    "2021-10-05 13:40:13.764576",
    "2021-10-05 13:40:13.764579",
    "2021-10-05 13:40:13.764582",
    "2021-10-05 13:40:13.764585",
    "2021-10-05 13:40:13.764588",
    "2021-10-05 13:40:13.764590",
    "2021-10-05 13:40:13.764593",
    "2021-10-05 13:40:13.764595",
    "2021-10-05 13:40:13.764598",
    "2021-10-05 13:40:13.764601",
    "2021-10-05 13:40:13.764603",
    "2021-10-05 13:40:13.764607",
    "2021-10-05 13:40:13.764609",
    "2021-10-05 13:40:13.764612",
    "2021-10-05 13:40:13.764614",
    "2021-10-05 13:40:13.764617",
    "2021-10-05 13:40:13.764620",
    "2021-10-05 13:40:13.764622",
    "2021-10-05 13:40:13.764625",
    "2021-10-05 13:40:13.764628",
    "2021-10-05 13:40:13.764631",
    "2021-10-05 13:40:13.764633",
    "2021-10-05 13:40:13.764636",
    "2021-10-05 13:40:13.764638",
    "2021-10-05 13:40:13.764641",
    "2021-10-05 13:40:13.764645",
    "2021-10-05 13:40:13.764650",
    "2021-10-05 13:40:13.764654",
    "2021-10-05 13:40:13.764659",
    "2021-10-05 13:40:13.764664",
    "2021-10-05 13:40:13.764669",
    "2021-10-05 13:40:13.764674",
    "2021-10-05 13:40:13.764679",
    "2021-10-05 13:40:13.764684",
    "2021-10-05 13:40:13.764689",
    "2021-10-05 13:40:13.764693",
    "2021-10-05 13:40:13.764697",
    "2021-10-05 13:40:13.764700",
    "2021-10-05 13:40:13.764702",
    "2021-10-05 13:40:13.764705",
    "2021-10-05 13:40:13.764707",
    "2021-10-05 13:40:13.764710",
    "2021-10-05 13:40:13.764713",
    "2021-10-05 13:40:13.764717",
    "2021-10-05 13:40:13.764720",
    "2021-10-05 13:40:13.764723",
    "2021-10-05 13:40:13.764726",
    "2021-10-05 13:40:13.764728",
    "2021-10-05 13:40:13.764731",
    "2021-10-05 13:40:13.764734",
    "2021-10-05 13:40:13.764737",
    "2021-10-05 13:40:13.764739",
    "2021-10-05 13:40:13.764742",
    "2021-10-05 13:40:13.764744",
    "2021-10-05 13:40:13.764747",
    "2021-10-05 13:40:13.764751",
    "2021-10-05 13:40:13.764754",
    "2021-10-05 13:40:13.764757",
    "2021-10-05 13:40:13.764759",
    "2021-10-05 13:40:13.764762",
    "2021-10-05 13:40:13.764765",
    "2021-10-05 13:40:13.764767",
    "2021-10-05 13:40:13.764770",
    "2021-10-05 13:40:13.764773",
    "2021-10-05 13:40:13.764776",
    "2021-10-05 13:40:13.764778",
    "2021-10-05 13:40:13.764781",
    "2021-10-05 13:40:13.764784",
    "2021-10-05 13:40:13.764786",
    "2021-10-05 13:40:13.764789",
    "2021-10-05 13:40:13.764791",
    "2021-10-05 13:40:13.764795",
    "2021-10-05 13:40:13.764798",
    "2021-10-05 13:40:13.764801",
    "2021-10-05 13:40:13.764804",
    "2021-10-05 13:40:13.764806",
    "2021-10-05 13:40:13.764809",
    "2021-10-05 13:40:13.764811",
    "2021-10-05 13:40:13.764814",
    "2021-10-05 13:40:13.764817",
    "2021-10-05 13:40:13.764819",
    "2021-10-05 13:40:13.764822",
    "2021-10-05 13:40:13.764825",
    "2021-10-05 13:40:13.764827",
    "2021-10-05 13:40:13.764830",
    "2021-10-05 13:40:13.764833",
    "2021-10-05 13:40:13.764836",
    "2021-10-05 13:40:13.764840",
    "2021-10-05 13:40:13.764843",
    "2021-10-05 13:40:13.764846",
    "2021-10-05 13:40:13.764849",
    "2021-10-05 13:40:13.764851",
    "2021-10-05 13:40:13.764854",
    "2021-10-05 13:40:13.764856",
    "2021-10-05 13:40:13.764859",
    "2021-10-05 13:40:13.764862",
    "2021-10-05 13:40:13.764864",
    "2021-10-05 13:40:13.764867",
    "2021-10-05 13:40:13.764869",
    "2021-10-05 13:40:13.764872",
    "2021-10-05 13:40:13.764875",
    "2021-10-05 13:40:13.764877",
    "2021-10-05 13:40:13.764880",
    "2021-10-05 13:40:13.764882",
    "2021-10-05 13:40:13.764885",
    "2021-10-05 13:40:13.764888",
    "2021-10-05 13:40:13.764890",
    "2021-10-05 13:40:13.764895",
    "2021-10-05 13:40:13.764898",
    "2021-10-05 13:40:13.764900",
    "2021-10-05 13:40:13.764903",
    "2021-10-05 13:40:13.764906",
    "2021-10-05 13:40:13.764909",
    "2021-10-05 13:40:13.764911",
    "2021-10-05 13:40:13.764914",
    "2021-10-05 13:40:13.764917",
    "2021-10-05 13:40:13.764919",
    "2021-10-05 13:40:13.764922",
    "2021-10-05 13:40:13.764925",
    "2021-10-05 13:40:13.764927",
    "2021-10-05 13:40:13.764930",
    "2021-10-05 13:40:13.764932",
    "2021-10-05 13:40:13.764935",
    "2021-10-05 13:40:13.764938",
    "2021-10-05 13:40:13.764940",
    "2021-10-05 13:40:13.764943",
    "2021-10-05 13:40:13.764946",
    "2021-10-05 13:40:13.764948",
    "2021-10-05 13:40:13.764951",
    "2021-10-05 13:40:13.764954",
    "2021-10-05 13:40:13.764957",
    "2021-10-05 13:40:13.764959",
    "2021-10-05 13:40:13.764962",
    "2021-10-05 13:40:13.764965",
    "2021-10-05 13:40:13.764967",
    "2021-10-05 13:40:13.764970",
    "2021-10-05 13:40:13.764972",
    "2021-10-05 13:40:13.764975",
    "2021-10-05 13:40:13.764977",
    "2021-10-05 13:40:13.764980",
    "2021-10-05 13:40:13.764983",
    "2021-10-05 13:40:13.764985",
    "2021-10-05 13:40:13.764988",
    "2021-10-05 13:40:13.764990",
    "2021-10-05 13:40:13.764993",
    "2021-10-05 13:40:13.764995",
    "2021-10-05 13:40:13.764998",
    "2021-10-05 13:40:13.765000",
    "2021-10-05 13:40:13.765003",
    "2021-10-05 13:40:13.765006",
    "2021-10-05 13:40:13.765008",
    "2021-10-05 13:40:13.765012",
    "2021-10-05 13:40:13.765016",
    "2021-10-05 13:40:13.765018",
    "2021-10-05 13:40:13.765021",
    "2021-10-05 13:40:13.765023",
    "2021-10-05 13:40:13.765026",
    "2021-10-05 13:40:13.765029",
    "2021-10-05 13:40:13.765031",
    "2021-10-05 13:40:13.765034",
    "2021-10-05 13:40:13.765037",
    "2021-10-05 13:40:13.765039",
    "2021-10-05 13:40:13.765042",
    "2021-10-05 13:40:13.765044",
    "2021-10-05 13:40:13.765047",
    "2021-10-05 13:40:13.765049",
    "2021-10-05 13:40:13.765052",
    "2021-10-05 13:40:13.765054",
    "2021-10-05 13:40:13.765057",
    "2021-10-05 13:40:13.765059",
    "2021-10-05 13:40:13.765061",
    "2021-10-05 13:40:13.765064",
    "2021-10-05 13:40:13.765067",
    "2021-10-05 13:40:13.765070",
    "2021-10-05 13:40:13.765072",
    "2021-10-05 13:40:13.765075",
    "2021-10-05 13:40:13.765077",
    "2021-10-05 13:40:13.765080",
    "2021-10-05 13:40:13.765082",
    "2021-10-05 13:40:13.765085",
    "2021-10-05 13:40:13.765088",
    "2021-10-05 13:40:13.765090",
    "2021-10-05 13:40:13.765093",
    "2021-10-05 13:40:13.765096",
    "2021-10-05 13:40:13.765098",
    "2021-10-05 13:40:13.765101",
    "2021-10-05 13:40:13.765103",
    "2021-10-05 13:40:13.765106",
    "2021-10-05 13:40:13.765109",
    "2021-10-05 13:40:13.765111",
    "2021-10-05 13:40:13.765113",
    "2021-10-05 13:40:13.765116",
    "2021-10-05 13:40:13.765119",
    "2021-10-05 13:40:13.765121",
    "2021-10-05 13:40:13.765124",
    "2021-10-05 13:40:13.765126",
    "2021-10-05 13:40:13.765129",
    "2021-10-05 13:40:13.765132",
    "2021-10-05 13:40:13.765134",
    "2021-10-05 13:40:13.765137",
    "2021-10-05 13:40:13.765139",
    "2021-10-05 13:40:13.765142",
    "2021-10-05 13:40:13.765144",
    "2021-10-05 13:40:13.765147",
    "2021-10-05 13:40:13.765149",
    "2021-10-05 13:40:13.765152",
    "2021-10-05 13:40:13.765154",
    "2021-10-05 13:40:13.765157",
    "2021-10-05 13:40:13.765159",
    "2021-10-05 13:40:13.765162",
    "2021-10-05 13:40:13.765164",
    "2021-10-05 13:40:13.765168",
    "2021-10-05 13:40:13.765172",
    "2021-10-05 13:40:13.765174",
    "2021-10-05 13:40:13.765177",
    "2021-10-05 13:40:13.765179",
    "2021-10-05 13:40:13.765182",
    "2021-10-05 13:40:13.765184",
    "2021-10-05 13:40:13.765187",
    "2021-10-05 13:40:13.765189",
    "2021-10-05 13:40:13.765192",
    "2021-10-05 13:40:13.765195",
    "2021-10-05 13:40:13.765198",
    "2021-10-05 13:40:13.765200",
    "2021-10-05 13:40:13.765203",
    "2021-10-05 13:40:13.765205",
    "2021-10-05 13:40:13.765208",
    "2021-10-05 13:40:13.765210",
    "2021-10-05 13:40:13.765213",
    "2021-10-05 13:40:13.765216",
    "2021-10-05 13:40:13.765218",
    "2021-10-05 13:40:13.765221",
    "2021-10-05 13:40:13.765223",
    "2021-10-05 13:40:13.765226",
    "2021-10-05 13:40:13.765228",
    "2021-10-05 13:40:13.765231",
    "2021-10-05 13:40:13.765233",
    "2021-10-05 13:40:13.765236",
    "2021-10-05 13:40:13.765238",
    "2021-10-05 13:40:13.765241",
    "2021-10-05 13:40:13.765243",
    "2021-10-05 13:40:13.765246",
    "2021-10-05 13:40:13.765248",
    "2021-10-05 13:40:13.765251",
    "2021-10-05 13:40:13.765253",
    "2021-10-05 13:40:13.765256",
    "2021-10-05 13:40:13.765258",
    "2021-10-05 13:40:13.765261",
    "2021-10-05 13:40:13.765264",
    "2021-10-05 13:40:13.765266",
    "2021-10-05 13:40:13.765269",
    "2021-10-05 13:40:13.765272",
    "2021-10-05 13:40:13.765274",
    "2021-10-05 13:40:13.765277",
    "2021-10-05 13:40:13.765280",
    "2021-10-05 13:40:13.765282",
    "2021-10-05 13:40:13.765285",
    "2021-10-05 13:40:13.765287",
    "2021-10-05 13:40:13.765290",
    "2021-10-05 13:40:13.765292",
    "2021-10-05 13:40:13.765295",
    "2021-10-05 13:40:13.765298",
    "2021-10-05 13:40:13.765300",
    "2021-10-05 13:40:13.765303",
    "2021-10-05 13:40:13.765306",
    "2021-10-05 13:40:13.765309",
    "2021-10-05 13:40:13.765311",
    "2021-10-05 13:40:13.765314",
    "2021-10-05 13:40:13.765317",
    "2021-10-05 13:40:13.765319",
    "2021-10-05 13:40:13.765322",
    "2021-10-05 13:40:13.765325",
    "2021-10-05 13:40:13.765327",
    "2021-10-05 13:40:13.765330",
    "2021-10-05 13:40:13.765333",
    "2021-10-05 13:40:13.765336",
    "2021-10-05 13:40:13.765339"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
